# Weekly price-sheet update: a new weekly record is inserted as row 53
# (pushing the existing rows 53-70 down to 54-71), matching the pattern of
# one new observation being added to the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 53..70 down to 54..71, opening up a blank row 53.
$ws.Rows.Item(53).Insert()

# Populate the newly opened row 53 with the new weekly record.
$ws.Cells.Item(53, 1).Value  = 11
$ws.Cells.Item(53, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(53, 3).Value  = "Bíobío"
$ws.Cells.Item(53, 4).Value  = 44463
$ws.Cells.Item(53, 5).Value  = 8
$ws.Cells.Item(53, 6).Value  = 100112043
$ws.Cells.Item(53, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(53, 8).Value  = "Sin especificar"
$ws.Cells.Item(53, 9).Value  = "Primera"
$ws.Cells.Item(53, 10).Value = 100
$ws.Cells.Item(53, 11).Value = 15000
$ws.Cells.Item(53, 12).Value = 15500
$ws.Cells.Item(53, 13).Value = 15250
$ws.Cells.Item(53, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(53, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value = 254
$ws.Cells.Item(53, 17).Value = 60
$ws.Cells.Item(53, 18).Value = "Hortaliza"
